$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "DC_TYPE"
$ws.Range("B1").Value = "SUB_TYPE"
$ws.Range("C1").Value = "CELL_TYPE"
$ws.Range("D1").Value = "TEST_CASE"
$ws.Range("E1").Value = "TEST_ID"
$ws.Range("F1").Value = "SCOPE"
$ws.Range("G1").Value = "PHASE"
$ws.Range("H1").Value = "CELLS"
$ws.Range("I1").Value = "DESCRIPTION"
$ws.Range("J1").Value = "REQUIREMENTS"
$ws.Range("K1").Value = "STEPS"
$ws.Range("L1").Value = "MULTI_DRIVEWAY"
$ws.Range("M1").Value = "DRIVEWAY_TYPE"
$ws.Range("N1").Value = "COMBINED_TEST"
$ws.Range("O1").Value = "IMAGE"
$ws.Range("P1").Value = "LAST_MODIFIED"
$ws.Range("Q1").Value = "MODIFIED_USER"
